# The workbook gained one new weekly price-report row. It was inserted as
# row 210 (pushing the former rows 210-312 down to 211-313), so every
# subsequent row keeps its original data but moves down by one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 210, shifting rows 210:312 down to 211:313.
$ws.Rows.Item(210).Insert()

# Populate the newly inserted row 210 with the new record's data.
$ws.Range("A210").Value = 11
$ws.Range("B210").Value = 'Vega Monumental Concepción'
$ws.Range("C210").Value = 'Bíobío'
$ws.Range("D210").Value = 44845
$ws.Range("E210").Value = 8
$ws.Range("F210").Value = 100112009
$ws.Range("G210").Value = 'Acelga'
$ws.Range("H210").Value = 'Sin especificar'
$ws.Range("I210").Value = 'Primera'
$ws.Range("J210").Value = 110
$ws.Range("K210").Value = 700
$ws.Range("L210").Value = 750
$ws.Range("M210").Value = 727
$ws.Range("N210").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O210").Value = 'Región de Ñuble'
$ws.Range("P210").Value = 727
$ws.Range("Q210").Value = 1
$ws.Range("R210").Value = 'Hortaliza'
